$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages")

# ------------------------------------------------------------------
# 1) Fix a typo in an existing shared string (row 69, column A):
#    "minestones-and-feedback" -> "milestones-and-feedback"
# ------------------------------------------------------------------
$ws.Range("A69").Value2 = "https://heatlabs.net/blog/milestones-and-feedback"

# ------------------------------------------------------------------
# 2) A new tracked page was added at the bottom of the table
#    ("wishlist-heatlabs-on-steam"). Its initial checker status
#    (columns B:E) starts out as "NOT NEEDED", which is the same
#    status the previous last data row (275) used to have. The
#    checker-status block for every row from 82 down to 275 shifts
#    down by one row to make room (column A, the URL itself, is
#    untouched for those rows).
# ------------------------------------------------------------------

# 2a) Seed row 276 with row 275's current formatting + values (this
#     gives row 276 the "NOT NEEDED" x4 status block and the correct
#     cell style, matching what a freshly appended tracked row looks
#     like).
$ws.Range("A275:E275").Copy()
$ws.Range("A276:E276").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A275:E275").Copy()
$ws.Range("A276:E276").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# 2b) Shift the checker-status columns (B:E) for rows 82-275 down by
#     one row in a single bulk operation.
$shiftedBlock = $ws.Range("B81:E274").Value2
$ws.Range("B82:E275").Value2 = $shiftedBlock

# 2c) Give the new row its real URL.
$newUrl = "https://heatlabs.net/blog/wishlist-heatlabs-on-steam"
$ws.Range("A276").Value2 = $newUrl

# 2d) Hyperlink the new row's URL cell, same as every other row in
#     the table.
$ws.Hyperlinks.Add($ws.Range("A276"), $newUrl) | Out-Null

Write-Host "Edit applied."
